$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @("IMX-USD", "TAO-USD", "MNT-USD")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$row = $lastRow + 1
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
